$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Human" column (D), shifting
# "Human" to column E, and turn the new column D into "Copilot" with its
# own set of answers.
$ws.Columns("D:D").Insert()

$ws.Range("D1").Value = "Copilot"

$ws.Range("D2").Value  = "A"
$ws.Range("D3").Value  = "D"
$ws.Range("D4").Value  = "B"
$ws.Range("D5").Value  = "A"
$ws.Range("D6").Value  = "D"
$ws.Range("D7").Value  = "E"
$ws.Range("D8").Value  = "C"
$ws.Range("D9").Value  = "A"
$ws.Range("D10").Value = "C"
$ws.Range("D11").Value = "A"
$ws.Range("D12").Value = "C"
$ws.Range("D13").Value = "A"
$ws.Range("D14").Value = "F"
$ws.Range("D15").Value = "C"
$ws.Range("D16").Value = "B"
$ws.Range("D17").Value = "B"
$ws.Range("D18").Value = "E"
$ws.Range("D19").Value = "D"
$ws.Range("D20").Value = "A"
$ws.Range("D21").Value = "D"
$ws.Range("D22").Value = "A"
$ws.Range("D23").Value = "A"
$ws.Range("D24").Value = "D"
$ws.Range("D25").Value = "B"
$ws.Range("D26").Value = "C"
$ws.Range("D27").Value = "D"
$ws.Range("D28").Value = "A"
$ws.Range("D29").Value = "E"
$ws.Range("D30").Value = "B"
$ws.Range("D31").Value = "D"
$ws.Range("D32").Value = "D"
$ws.Range("D33").Value = "E"
$ws.Range("D34").Value = "B"
$ws.Range("D35").Value = "A"
$ws.Range("D36").Value = "A"
$ws.Range("D37").Value = "C"
$ws.Range("D38").Value = "D"
$ws.Range("D39").Value = "B"
$ws.Range("D40").Value = "E"
$ws.Range("D41").Value = "B"

# Match the author's final view/selection state.
$null = $ws.Range("F35").Select()
